# Apply crypto price/volume updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "61.650.67"
Set-TextCell $ws.Range("E2") "  +0.86%  "

Set-TextCell $ws.Range("D3") "3.393.75"
Set-TextCell $ws.Range("E3") "  -0.07%  "

Set-TextCell $ws.Range("E4") "  +0.04%  "

Set-TextCell $ws.Range("D5") "576.46"
Set-TextCell $ws.Range("E5") "  +0.56%  "

Set-TextCell $ws.Range("D6") "141.68"
Set-TextCell $ws.Range("E6") "  -0.35%  "

Set-TextCell $ws.Range("E7") "  +0.01%  "

Set-TextCell $ws.Range("D8") "0.474"
Set-TextCell $ws.Range("E8") "  -0.45%  "

Set-TextCell $ws.Range("D9") "7.68"
Set-TextCell $ws.Range("E9") "  +0.48%  "

Set-TextCell $ws.Range("E10") "  -1.13%  "

Set-TextCell $ws.Range("D11") "0.386"
Set-TextCell $ws.Range("E11") "  -2.30%  "

Set-TextCell $ws.Range("D12") "3.973.16"
Set-TextCell $ws.Range("E12") "  +0.00%  "

Set-TextCell $ws.Range("E13") "  -0.16%  "

Set-TextCell $ws.Range("E14") "  +1.54%  "

Set-TextCell $ws.Range("D15") "3.386.70"
Set-TextCell $ws.Range("E15") "  -0.25%  "

Set-TextCell $ws.Range("E16") "  -0.51%  "

Set-TextCell $ws.Range("D17") "61.695.06"
Set-TextCell $ws.Range("E17") "  +0.93%  "

Set-TextCell $ws.Range("D18") "6.13"
Set-TextCell $ws.Range("E18") "  +0.08%  "

Set-TextCell $ws.Range("D19") "13.63"
Set-TextCell $ws.Range("E19") "  -1.10%  "

Set-TextCell $ws.Range("D20") "9.01"
Set-TextCell $ws.Range("E20") "  +1.16%  "

Set-TextCell $ws.Range("D21") "391.17"
Set-TextCell $ws.Range("E21") "  +2.25%  "

Set-TextCell $ws.Range("D22") "74.88"
Set-TextCell $ws.Range("E22") "  +0.31%  "

Set-TextCell $ws.Range("D23") "0.549"
Set-TextCell $ws.Range("E23") "  -1.29%  "

Set-TextCell $ws.Range("E24") "  -0.16%  "

Set-TextCell $ws.Range("B25") "PEPE"
Set-TextCell $ws.Range("C25") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell $ws.Range("D25") "0.0000113"
Set-TextCell $ws.Range("E25") "  -3.64%  "

Set-TextCell $ws.Range("B26") "Kaspa"
Set-TextCell $ws.Range("C26") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws.Range("D26") "0.198"
Set-TextCell $ws.Range("E26") "  +9.03%  "

Set-TextCell $ws.Range("D27") "0.999"
Set-TextCell $ws.Range("E27") "  -0.05%  "

Set-TextCell $ws.Range("D28") "7.40"
Set-TextCell $ws.Range("E28") "  +0.49%  "

Set-TextCell $ws.Range("D29") "8.00"
Set-TextCell $ws.Range("E29") "  +0.04%  "

Set-TextCell $ws.Range("E30") "  -0.63%  "

Set-TextCell $ws.Range("D31") "1.43"
Set-TextCell $ws.Range("E31") "  +1.89%  "

Set-TextCell $ws.Range("E32") "  +0.00%  "

Set-TextCell $ws.Range("D33") "23.33"
Set-TextCell $ws.Range("E33") "  -0.60%  "

Set-TextCell $ws.Range("D34") "6.93"
Set-TextCell $ws.Range("E34") "  -0.96%  "

Set-TextCell $ws.Range("D35") "169.14"
Set-TextCell $ws.Range("E35") "  +1.42%  "

Set-TextCell $ws.Range("D36") "5.04"
Set-TextCell $ws.Range("E36") "  +0.63%  "

Set-TextCell $ws.Range("D37") "3.428.83"
Set-TextCell $ws.Range("E37") "  +0.03%  "

Set-TextCell $ws.Range("D38") "1.48"
Set-TextCell $ws.Range("E38") "  -0.80%  "

Set-TextCell $ws.Range("D39") "0.0765"
Set-TextCell $ws.Range("E39") "  -0.82%  "

Set-TextCell $ws.Range("D40") "25.80"
Set-TextCell $ws.Range("E40") "  -5.33%  "

Set-TextCell $ws.Range("E41") "  -0.31%  "

Set-TextCell $ws.Range("E42") "  +0.20%  "

Set-TextCell $ws.Range("D43") "1.65"
Set-TextCell $ws.Range("E43") "  -0.96%  "

Set-TextCell $ws.Range("E44") "  +1.64%  "

Set-TextCell $ws.Range("D45") "2.456.47"
Set-TextCell $ws.Range("E45") "  -0.65%  "

Set-TextCell $ws.Range("D46") "22.75"
Set-TextCell $ws.Range("E46") "  -0.86%  "

Set-TextCell $ws.Range("E47") "  -1.96%  "

Set-TextCell $ws.Range("E48") "  +0.09%  "

Set-TextCell $ws.Range("E49") "  -1.22%  "

Set-TextCell $ws.Range("E50") "  -5.86%  "

Set-TextCell $ws.Range("E51") "  -1.29%  "
